# "Add all field search" - update several sample cells in Sheet1 to new
# demo text (replacing long lorem-ipsum filler with short search-friendly
# strings), then move the viewport/selection to reflect where the user is
# now working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- replace long filler text / typo values with the new sample content --
$ws.Range("E6").Value  = "Brown fox brown dog"
$ws.Range("E7").Value  = "The quick brown fox jumps over the lazy dog"
$ws.Range("E8").Value  = "The quick brown fox jumps over the quick dog"
$ws.Range("D9").Value  = "Tzoli"
$ws.Range("E9").Value  = "dffefewfew"
$ws.Range("E10").Value = "Tzoli"

# row heights were tied to the old wrapped text; re-fit them now that the
# new, shorter strings live in rows 6-10
$ws.Rows("6:10").AutoFit()

# reflect the new scroll/selection position used while making this edit
[void]$ws.Range("G8").Select()
